# The data rows for row 19 and row 20 were swapped (their observation
# records traded places). Columns C, I, P, S, T, U, V, W, Y, Z, AA, AB,
# AD, AE, AG, AT, AW, AX, AY already hold identical values in both rows,
# so only A, B, D, E, F, G, H, Q, R actually need to be exchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $cell19 = $ws.Range($col + "19")
    $cell20 = $ws.Range($col + "20")

    # Use Value2 for reading to avoid locale/format based text conversion
    # and get the exact underlying number/string.
    $v19 = $cell19.Value2
    $v20 = $cell20.Value2

    $cell19.Value = $v20
    $cell20.Value = $v19
}
